$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '72.204.78'
$ws.Range('E2').Value = '  +0.67%  '
$ws.Range('D3').Value = '2.702.38'
$ws.Range('E3').Value = '  +2.96%  '
$ws.Range('E4').Value = '  +0.11%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '598.77'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  -1.19%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '174.77'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  -2.55%  '
$ws.Range('E7').Value = '  -0.01%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.524'
$ws.Range('D8').Style = "Normal"
$ws.Range('E8').Value = '  +0.01%  '
$ws.Range('D9').Value = '2.708.51'
$ws.Range('E9').Value = '  +3.18%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.167'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  -0.57%  '
$ws.Range('E11').Value = '  +2.03%  '
$ws.Range('E12').Value = '  +2.06%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '5.00'
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').Value = '  -0.52%  '
$ws.Range('D14').Value = '3.204.92'
$ws.Range('E14').Value = '  +3.92%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '0.0000185'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  -1.18%  '
$ws.Range('D16').Value = '71.970.13'
$ws.Range('E16').Value = '  +0.33%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '26.26'
$ws.Range('D17').Style = "Normal"
$ws.Range('E17').Value = '  -0.68%  '
$ws.Range('D18').Value = '2.708.95'
$ws.Range('E18').Value = '  +3.66%  '
$ws.Range('E19').Value = '  +6.97%  '
$ws.Range('E20').Value = '  +0.25%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '372.87'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  -2.45%  '
$ws.Range('E22').Value = '  +0.81%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '2.00'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  +0.47%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '72.33'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  -0.51%  '
$ws.Range('E25').Value = '  -0.20%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '4.36'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  -1.90%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '9.79'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  -1.71%  '
$ws.Range('D28').Value = '2.848.50'
$ws.Range('E28').Value = '  +3.12%  '
$ws.Range('E29').Value = '  -0.01%  '
$ws.Range('D30').Value = '0.0₃0977'
$ws.Range('E30').Value = '  +1.78%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '8.11'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  +0.80%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '504.10'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  -7.74%  '
$ws.Range('E33').Value = '  -1.33%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '1.82'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '  -0.35%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '1.00'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  +0.24%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '163.30'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  -2.01%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '19.66'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  +2.46%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '19.13'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  +0.02%  '
$ws.Range('E39').Value = '  -0.24%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.107'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  -6.07%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '1.79'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  -3.31%  '
$ws.Range('B42').Value = 'USDe'
$ws.Range('C42').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '1.00'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  -0.02%  '
$ws.Range('B43').Value = 'RenderToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '5.05'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  +0.52%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '0.334'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  +0.64%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '2.55'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  -2.45%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '157.68'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  +4.67%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '39.53'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  +0.42%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '0.564'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  +4.78%  '
$ws.Range('E49').Value = '  +3.48%  '
$ws.Range('E50').Value = '  +5.36%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '0.0763'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  +1.25%  '
